# Update the graphSpec sheet to refresh the x-axis labels/ranges based on
# data from the excel file (visScatterPlotter x-axis update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("graphSpec")

# Row 3: Ego Speed plot
$ws.Range("B3").Value = "Ego Speed"
$ws.Range("C3").Value = "vehSpd"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 85
$ws.Range("F3").Value = "Ego Speed (km/h)"

# Row 4: Dead Time plot - x axis label
$ws.Range("F4").Value = "Time (sec)"

# Row 5: Intervention Duration plot - x axis label
$ws.Range("F5").Value = "Time (sec)"

# Row 6: Max Throttle Value plot - x axis label
$ws.Range("F6").Value = "Percentage (%)"

# Row 7: Max Steering Angle plot - x axis label
$ws.Range("F7").Value = "Angle (rad)"

# Row 8: Max Steering Angle Rate plot - x axis label
$ws.Range("F8").Value = "Angle Speed (rad/s)"

# Row 9: Max Yaw Rate plot - x axis label
$ws.Range("F9").Value = "Angle Speed (rad/s)"

# Row 10: Max Lateral Acceleration plot - fix typo + x axis label
$ws.Range("B10").Value = "Max Lateral Acceleration [m/s2]"
$ws.Range("F10").Value = "Acceleration (m/s2)"

# Row 11: Longitudinal Clearance plot - x axis label
$ws.Range("F11").Value = "Distance (m)"

# Update the sheet view so the previously selected cell scrolls into place.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("E19").Select()
